$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns for the affected rows remain text, matching original inlineStr cells
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.105.20"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "1.560.08"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "288.37"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "0.3791"
$ws.Range("E7").Value = "  +2.65%  "

$ws.Range("D8").Value = "0.3273"
$ws.Range("E8").Value = "  -1.71%  "

$ws.Range("D9").Value = "43.27"
$ws.Range("E9").Value = "  -9.56%  "

$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").Value = "0.07345"
$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "19.87"
$ws.Range("E13").Value = "  -4.64%  "

$ws.Range("D14").Value = "5.794"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").Value = "6.875"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "1.561.02"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("D18").Value = "0.06652"

$ws.Range("D19").Value = "85.58"
$ws.Range("E19").Value = "  -3.08%  "

$ws.Range("D20").Value = "6.449"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "16.07"
$ws.Range("E22").Value = "  -2.78%  "

$ws.Range("D23").Value = "11.66"
$ws.Range("E23").Value = "  -3.16%  "

$ws.Range("D24").Value = "22.127.36"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").Value = "2.261"
$ws.Range("E25").Value = "  -5.23%  "

$ws.Range("D26").Value = "2.537"
$ws.Range("E26").Value = "  -3.91%  "

$ws.Range("D27").Value = "150.02"
$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").Value = "4.867"
$ws.Range("E29").Value = "  -2.43%  "

$ws.Range("D30").Value = "1.737.80"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").Value = "121.07"
$ws.Range("E31").Value = "  -3.58%  "

$ws.Range("D32").Value = "1.114"
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("D33").Value = "5.998"
$ws.Range("E33").Value = "  -1.89%  "

$ws.Range("D34").Value = "1.802"
$ws.Range("E34").Value = "  -9.77%  "

$ws.Range("D35").Value = "9.319"
$ws.Range("E35").Value = "  -5.57%  "

$ws.Range("D36").Value = "0.08154"
$ws.Range("E36").Value = "  -2.72%  "

$ws.Range("D37").Value = "5.258"
$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("D38").Value = "0.06182"
$ws.Range("E38").Value = "  -3.36%  "

$ws.Range("D39").Value = "0.02289"
$ws.Range("E39").Value = "  -7.13%  "

$ws.Range("D40").Value = "0.2137"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("D41").Value = "1.226"
$ws.Range("E41").Value = "  -5.57%  "

$ws.Range("D42").Value = "10.99"
$ws.Range("E42").Value = "  -4.25%  "

$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("D44").Value = "0.5969"
$ws.Range("E44").Value = "  -4.93%  "

$ws.Range("D45").Value = "13.70"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("D46").Value = "3.739"
$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("D47").Value = "0.5766"
$ws.Range("E47").Value = "  -5.62%  "

$ws.Range("D48").Value = "1.966"
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("D49").Value = "120.18"
$ws.Range("E49").Value = "  -4.14%  "

$ws.Range("D50").Value = "1.169"
$ws.Range("E50").Value = "  -3.57%  "

$ws.Range("D51").Value = "0.06971"
$ws.Range("E51").Value = "  -3.49%  "
